# Commit: "update 038 update 012"
# Updates the point-template coordinates for templates "012" (row 13) and
# "038" (row 39) on the "Templates" sheet, and fixes row 39's File Name
# column (F39) which incorrectly pointed at template "012" instead of "038".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templates")
$ws.Activate()

# --- Template "012" (row 13): new X/Y start/end coordinates ---
$ws.Range("B13").Value = 858
$ws.Range("C13").Value = 954
$ws.Range("D13").Value = 1007
$ws.Range("E13").Value = 1000

# --- Template "038" (row 39): new X/Y start/end coordinates, and fix the
#     File Name column which was still referencing "012" ---
$ws.Range("B39").Value = 858
$ws.Range("C39").Value = 954
$ws.Range("D39").Value = 1007
$ws.Range("E39").Value = 1000
$ws.Range("F39").Value = "038"

# --- Reflect the selection left behind after making the edit ---
$ws.Range("E39").Select() | Out-Null
